$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the current row 5 ("Stadion"), before the
# current row 6 ("Sprungschanze (Anlauf)"), to hold the two new
# "Stadion, ..." values.
$ws.Rows("6:7").Insert()

$ws.Range("A6").Value = "Stadion, überdacht"
$ws.Range("B6").Value = 1441
$ws.Range("C6").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

$ws.Range("A7").Value = "Stadion, nicht überdacht"
$ws.Range("B7").Value = 1442
$ws.Range("C7").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

# Append a new row 10 for "Wassersportanlage" after the existing rows
# (now shifted to 8 and 9).
$ws.Range("A10").Value = "Wassersportanlage"
$ws.Range("B10").Value = 1650
$ws.Range("C10").Value = "http://inspire.ec.europa.eu/codelist/CurrentUseValue/publicServices"

# New row 10 falls outside the previously formatted range, so copy the
# "link-like" cell style used by column C from a neighboring row.
$ws.Range("C10").Style = $ws.Range("C9").Style

# Remove the hyperlink that used to span C3:C7 (now effectively C3:C9)
# entirely, leaving only the single hyperlink on C2 untouched.
$ws.Hyperlinks.Item(2).Delete()

# Leave the cursor on the newly added last row, mirroring where the
# author's selection ended up after the edit.
$ws.Range("C10").Select() | Out-Null
